$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.254.22"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "1.866.36"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "243.66"
$ws.Range("E5").Value = "  +3.76%  "

# Row 6
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  +0.67%  "

# Row 8
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "42.65"
$ws.Range("E8").Value = "  -2.33%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2867"
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("D10").Value = "0.06478"
$ws.Range("E10").Value = "  -1.45%  "

# Row 11
$ws.Range("D11").Value = "21.11"
$ws.Range("E11").Value = "  -2.49%  "

# Row 12
$ws.Range("D12").Value = "0.07733"
$ws.Range("E12").Value = "  -1.96%  "

# Row 13
$ws.Range("D13").Value = "1.885.77"
$ws.Range("E13").Value = "  +0.94%  "

# Row 14
$ws.Range("D14").Value = "95.08"
$ws.Range("E14").Value = "  -1.53%  "

# Row 15
$ws.Range("D15").Value = "0.7060"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("D16").Value = "5.104"
$ws.Range("E16").Value = "  -0.11%  "

# Row 17
$ws.Range("D17").Value = "275.57"
$ws.Range("E17").Value = "  +2.73%  "

# Row 18
$ws.Range("D18").Value = "30.238.10"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").Value = "  -4.69%  "

# Row 20
$ws.Range("D20").Value = "0.000007562"
$ws.Range("E20").Value = "  -1.54%  "

# Row 21
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "2.120.81"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").Value = "5.212"
$ws.Range("E24").Value = "  -0.50%  "

# Row 25
$ws.Range("D25").Value = "6.127"
$ws.Range("E25").Value = "  -1.06%  "

# Row 26
$ws.Range("D26").Value = "9.301"
$ws.Range("E26").Value = "  -1.08%  "

# Row 27
$ws.Range("D27").Value = "165.01"
$ws.Range("E27").Value = "  -1.52%  "

# Row 28
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$ws.Range("D29").Value = "1.906"
$ws.Range("E29").Value = "  -2.14%  "

# Row 30
$ws.Range("D30").Value = "1.372"
$ws.Range("E30").Value = "  +0.83%  "

# Row 31
$ws.Range("D31").Value = "0.09860"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").Value = "1.512"
$ws.Range("E32").Value = "  +3.55%  "

# Row 33
$ws.Range("D33").Value = "4.249"
$ws.Range("E33").Value = "  -3.37%  "

# Row 34
$ws.Range("D34").Value = "4.027"
$ws.Range("E34").Value = "  -1.13%  "

# Row 35
$ws.Range("D35").Value = "0.04730"
$ws.Range("E35").Value = "  -0.66%  "

# Row 36
$ws.Range("D36").Value = "1.119"
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
$ws.Range("D37").Value = "0.6920"
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("D38").Value = "2.706"
$ws.Range("E38").Value = "  -0.66%  "

# Row 39
$ws.Range("D39").Value = "0.01842"
$ws.Range("E39").Value = "  -1.70%  "

# Row 40
$ws.Range("D40").Value = "2.738"
$ws.Range("E40").Value = "  -2.29%  "

# Row 41
$ws.Range("D41").Value = "6.332"
$ws.Range("E41").Value = "  +1.52%  "

# Row 42
$ws.Range("D42").Value = "70.17"
$ws.Range("E42").Value = "  -4.42%  "

# Row 43
$ws.Range("D43").Value = "0.8413"
$ws.Range("E43").Value = "  -0.10%  "

# Row 44
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "1.896"
$ws.Range("E45").Value = "  -2.96%  "

# Row 46
$ws.Range("D46").Value = "0.4084"
$ws.Range("E46").Value = "  -2.21%  "

# Row 47
$ws.Range("D47").Value = "101.91"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("D48").Value = "9.263"
$ws.Range("E48").Value = "  +1.72%  "

# Row 49
$ws.Range("D49").Value = "7.063"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("D50").Value = "921.51"
$ws.Range("E50").Value = "  -4.72%  "

# Row 51
$ws.Range("D51").Value = "34.96"
$ws.Range("E51").Value = "  +1.26%  "

